$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.416.18"
$ws.Range("E2").Value = "  -3.03%  "
$ws.Range("D3").Value = "2.488.09"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'313.77"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "'94.78"
$ws.Range("E6").Value = "  -5.41%  "
$ws.Range("D7").Value = "'0.551"
$ws.Range("E7").Value = "  -2.54%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -4.03%  "
$ws.Range("D10").Value = "'33.65"
$ws.Range("E10").Value = "  -5.04%  "
$ws.Range("D11").Value = "'0.0784"
$ws.Range("E11").Value = "  -2.41%  "
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "'7.02"
$ws.Range("E13").Value = "  -3.56%  "
$ws.Range("D14").Value = "2.866.86"
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").Value = "'15.50"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "2.448.63"
$ws.Range("E16").Value = "  -7.36%  "
$ws.Range("D17").Value = "'0.796"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").Value = "41.415.01"
$ws.Range("E18").Value = "  -3.03%  "
$ws.Range("D19").Value = "'6.34"
$ws.Range("E19").Value = "  -5.26%  "
$ws.Range("D20").Value = "0.0₃0930"
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("D21").Value = "'11.26"
$ws.Range("E21").Value = "  -8.39%  "
$ws.Range("D22").Value = "'68.93"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").Value = "'237.81"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("E24").Value = "  -3.40%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "'1.91"
$ws.Range("E25").Value = "  -4.83%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'24.23"
$ws.Range("E27").Value = "  -5.33%  "
$ws.Range("D28").Value = "'2.23"
$ws.Range("E28").Value = "  -4.42%  "
$ws.Range("D29").Value = "'9.74"
$ws.Range("E29").Value = "  -3.66%  "
$ws.Range("D30").Value = "'36.52"
$ws.Range("E30").Value = "  -4.48%  "
$ws.Range("D31").Value = "'152.41"
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("D32").Value = "'5.49"
$ws.Range("E32").Value = "  -6.70%  "
$ws.Range("E33").Value = "  -3.10%  "
$ws.Range("D34").Value = "'2.56"
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("D35").Value = "'0.0750"
$ws.Range("E35").Value = "  -5.07%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").Value = "'17.85"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'3.09"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").Value = "'1.89"
$ws.Range("E38").Value = "  -3.90%  "
$ws.Range("D39").Value = "'0.115"
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("D40").Value = "'0.102"
$ws.Range("E40").Value = "  -7.75%  "
$ws.Range("E41").Value = "  +2.41%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "2.008.97"
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("D44").Value = "'19.56"
$ws.Range("E44").Value = "  -10.53%  "
$ws.Range("D45").Value = "'0.0287"
$ws.Range("E45").Value = "  -3.61%  "
$ws.Range("D46").Value = "'3.03"
$ws.Range("E46").Value = "  -8.95%  "
$ws.Range("D47").Value = "'8.78"
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("D48").Value = "2.726.93"
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("D49").Value = "'69.92"
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("D50").Value = "'97.64"
$ws.Range("E50").Value = "  -3.53%  "
$ws.Range("D51").Value = "'0.179"
$ws.Range("E51").Value = "  -6.05%  "
